$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A (GioiTinh) : new rows 3-4 keep "Cai" ---
$ws.Range("A3").Value = "Cái"
$ws.Range("A4").Value = "Cái"

# --- C (TrongLuong) : row2 corrected, rows 3-4 new ---
$ws.Range("C2").Value = 32
$ws.Range("C3").Value = 31
$ws.Range("C4").Value = 30

# --- D (TenLoai) : "Heo con" -> "Heo thit" for row2, filled down for 3-4 ---
$ws.Range("D2").Value = "Heo thịt"
$ws.Range("D3").Value = "Heo thịt"
$ws.Range("D4").Value = "Heo thịt"

# --- E (TenGiong) : row2 unchanged, rows 3-4 new breeds ---
$ws.Range("E3").Value = "Heo Móng Cái"
$ws.Range("E4").Value = "Heo Yorkshire"

# --- H (MaChuong) : row2 unchanged, rows 3-4 new pen codes ---
$ws.Range("H3").Value = "NT002"
$ws.Range("H4").Value = "NT003"

# --- B (NgaySinh) : "12/12/2022" -> "1/12/2022" (kept as text) ---
$ws.Range("B2").Value = "'1/12/2022"
$ws.Range("B3").Value = "'1/12/2022"
$ws.Range("B4").Value = "'1/12/2022"

# --- I (TinhTrang) : spelling fix "khoe" -> "khoe" (correct diacritic) ---
$ws.Range("I2").Value = "Sức khoẻ tốt"
$ws.Range("I3").Value = "Sức khoẻ tốt"
$ws.Range("I4").Value = "Sức khoẻ tốt"

# --- J (NguonGoc) : unchanged value filled down for new rows ---
$ws.Range("J3").Value = "Nhập ngoài"
$ws.Range("J4").Value = "Nhập ngoài"

# --- Selection moves to I8 ---
$null = $ws.Range("I8").Select()
